# Auto-generated Excel COM-interop script
# Applies updated market-price data (columns H-N) scraped by the scheduled pricing runner
# to rows across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row11
$ws.Cells.Item(11, 8).Value = 4116.6943
$ws.Cells.Item(11, 9).Value = 4116.6943
$ws.Cells.Item(11, 11).Value = 4116.6943
$ws.Cells.Item(11, 13).Value = -3976.6943

# ALC!row17
$ws.Cells.Item(17, 8).Value = 2241.3333
$ws.Cells.Item(17, 10).Value = 2241.3333
$ws.Cells.Item(17, 12).Value = 6723.999899999999
$ws.Cells.Item(17, 14).Value = -7059.999899999999

# ALC!row86
$ws.Cells.Item(86, 8).Value = 2954.111
$ws.Cells.Item(86, 9).Value = 1499.75
$ws.Cells.Item(86, 10).Value = 4117.6
$ws.Cells.Item(86, 11).Value = 1499.75
$ws.Cells.Item(86, 12).Value = 4117.6
$ws.Cells.Item(86, 13).Value = -376.75
$ws.Cells.Item(86, 14).Value = -6363.6

# ALC!row89
$ws.Cells.Item(89, 8).Value = 2954.111
$ws.Cells.Item(89, 9).Value = 1499.75
$ws.Cells.Item(89, 10).Value = 4117.6
$ws.Cells.Item(89, 11).Value = 7498.75
$ws.Cells.Item(89, 12).Value = 20588
$ws.Cells.Item(89, 13).Value = -1882.75
$ws.Cells.Item(89, 14).Value = -31820

# ALC!row98
$ws.Cells.Item(98, 8).Value = 1950.4
$ws.Cells.Item(98, 9).Value = 1824.5
$ws.Cells.Item(98, 11).Value = 1824.5
$ws.Cells.Item(98, 13).Value = -326.5

# ALC!row112
$ws.Cells.Item(112, 8).Value = 2635.8125
$ws.Cells.Item(112, 10).Value = 2635.8125
$ws.Cells.Item(112, 12).Value = 7907.4375
$ws.Cells.Item(112, 14).Value = -10123.4375

# ALC!row122
$ws.Cells.Item(122, 8).Value = 1950.4
$ws.Cells.Item(122, 9).Value = 1824.5
$ws.Cells.Item(122, 11).Value = 5473.5
$ws.Cells.Item(122, 13).Value = -3023.5

# ALC!row135
$ws.Cells.Item(135, 8).Value = 842.7241
$ws.Cells.Item(135, 9).Value = 794.5
$ws.Cells.Item(135, 11).Value = 7150.5
$ws.Cells.Item(135, 13).Value = -4615.5

# ALC!row137
$ws.Cells.Item(137, 8).Value = 4363.085
$ws.Cells.Item(137, 9).Value = 2617.0264
$ws.Cells.Item(137, 11).Value = 7851.0792
$ws.Cells.Item(137, 13).Value = -5301.0792

$ws = $wb.Worksheets.Item("ARM")
# ARM!row26
$ws.Cells.Item(26, 8).Value = 5106.2856
$ws.Cells.Item(26, 9).Value = 1029.2
$ws.Cells.Item(26, 10).Value = 15299
$ws.Cells.Item(26, 11).Value = 1029.2
$ws.Cells.Item(26, 12).Value = 15299
$ws.Cells.Item(26, 13).Value = -699.2
$ws.Cells.Item(26, 14).Value = -15959

# ARM!row45
$ws.Cells.Item(45, 8).Value = 8453.267
$ws.Cells.Item(45, 9).Value = 34299.668
$ws.Cells.Item(45, 11).Value = 34299.668
$ws.Cells.Item(45, 13).Value = -33922.668

# ARM!row61
$ws.Cells.Item(61, 8).Value = 3324.2344
$ws.Cells.Item(61, 10).Value = 8484.166999999999
$ws.Cells.Item(61, 12).Value = 8484.166999999999
$ws.Cells.Item(61, 14).Value = -8908.166999999999

# ARM!row74
$ws.Cells.Item(74, 8).Value = 288275.62
$ws.Cells.Item(74, 9).Value = 418058.53
$ws.Cells.Item(74, 11).Value = 418058.53
$ws.Cells.Item(74, 13).Value = -417184.53

# ARM!row77
$ws.Cells.Item(77, 8).Value = 288275.62
$ws.Cells.Item(77, 9).Value = 418058.53
$ws.Cells.Item(77, 11).Value = 2090292.65
$ws.Cells.Item(77, 13).Value = -2085924.65

# ARM!row97
$ws.Cells.Item(97, 8).Value = 2378045.2
$ws.Cells.Item(97, 9).Value = 3094109.5
$ws.Cells.Item(97, 11).Value = 3094109.5
$ws.Cells.Item(97, 13).Value = -3093613.5

# ARM!row132
$ws.Cells.Item(132, 8).Value = 4113.7886
$ws.Cells.Item(132, 9).Value = 3324.8572
$ws.Cells.Item(132, 11).Value = 9974.571599999999
$ws.Cells.Item(132, 13).Value = -7444.571599999999

# ARM!row133
$ws.Cells.Item(133, 8).Value = 95065
$ws.Cells.Item(133, 9).Value = 20000
$ws.Cells.Item(133, 10).Value = 120086.664
$ws.Cells.Item(133, 11).Value = 20000
$ws.Cells.Item(133, 12).Value = 120086.664
$ws.Cells.Item(133, 13).Value = -17470
$ws.Cells.Item(133, 14).Value = -125146.664

# ARM!row136
$ws.Cells.Item(136, 8).Value = 3324.2344
$ws.Cells.Item(136, 10).Value = 8484.166999999999
$ws.Cells.Item(136, 12).Value = 25452.501
$ws.Cells.Item(136, 14).Value = -30552.501

$ws = $wb.Worksheets.Item("BSM")
# BSM!row86
$ws.Cells.Item(86, 8).Value = 81481.24000000001
$ws.Cells.Item(86, 9).Value = 1407.7059
$ws.Cells.Item(86, 11).Value = 1407.7059
$ws.Cells.Item(86, 13).Value = -284.7058999999999

# BSM!row89
$ws.Cells.Item(89, 8).Value = 81481.24000000001
$ws.Cells.Item(89, 9).Value = 1407.7059
$ws.Cells.Item(89, 11).Value = 7038.5295
$ws.Cells.Item(89, 13).Value = -1422.5295

# BSM!row107
$ws.Cells.Item(107, 8).Value = 1508.6428
$ws.Cells.Item(107, 9).Value = 1600.3334
$ws.Cells.Item(107, 11).Value = 1600.3334
$ws.Cells.Item(107, 13).Value = 319.6666

$ws = $wb.Worksheets.Item("CRP")
# CRP!row31
$ws.Cells.Item(31, 8).Value = 45460384
$ws.Cells.Item(31, 9).Value = 142858750
$ws.Cells.Item(31, 10).Value = 7815.2666
$ws.Cells.Item(31, 11).Value = 142858750
$ws.Cells.Item(31, 12).Value = 7815.2666
$ws.Cells.Item(31, 13).Value = -142858455
$ws.Cells.Item(31, 14).Value = -8405.266599999999

# CRP!row34
$ws.Cells.Item(34, 8).Value = 45460384
$ws.Cells.Item(34, 9).Value = 142858750
$ws.Cells.Item(34, 10).Value = 7815.2666
$ws.Cells.Item(34, 11).Value = 142858750
$ws.Cells.Item(34, 12).Value = 7815.2666
$ws.Cells.Item(34, 13).Value = -142858548
$ws.Cells.Item(34, 14).Value = -8219.266599999999

# CRP!row58
$ws.Cells.Item(58, 8).Value = 3691.4285
$ws.Cells.Item(58, 9).Value = 2726.9473
$ws.Cells.Item(58, 11).Value = 2726.9473
$ws.Cells.Item(58, 13).Value = -2523.9473

# CRP!row74
$ws.Cells.Item(74, 8).Value = 40805.062
$ws.Cells.Item(74, 10).Value = 40805.062
$ws.Cells.Item(74, 12).Value = 40805.062
$ws.Cells.Item(74, 14).Value = -42553.062

# CRP!row77
$ws.Cells.Item(77, 8).Value = 40805.062
$ws.Cells.Item(77, 10).Value = 40805.062
$ws.Cells.Item(77, 12).Value = 122415.186
$ws.Cells.Item(77, 14).Value = -131151.186

# CRP!row105
$ws.Cells.Item(105, 8).Value = 789.7646999999999
$ws.Cells.Item(105, 9).Value = 782.875
$ws.Cells.Item(105, 10).Value = 900
$ws.Cells.Item(105, 11).Value = 782.875
$ws.Cells.Item(105, 12).Value = 900
$ws.Cells.Item(105, 13).Value = 964.125
$ws.Cells.Item(105, 14).Value = -4394

# CRP!row107
$ws.Cells.Item(107, 8).Value = 3500
$ws.Cells.Item(107, 10).Value = 5000
$ws.Cells.Item(107, 12).Value = 5000
$ws.Cells.Item(107, 14).Value = -8840

# CRP!row132
$ws.Cells.Item(132, 8).Value = 4611.5415
$ws.Cells.Item(132, 9).Value = 3327.4595
$ws.Cells.Item(132, 10).Value = 8930.727999999999
$ws.Cells.Item(132, 11).Value = 9982.378499999999
$ws.Cells.Item(132, 12).Value = 26792.184
$ws.Cells.Item(132, 13).Value = -7452.378499999999
$ws.Cells.Item(132, 14).Value = -31852.184

# CRP!row133
$ws.Cells.Item(133, 8).Value = 68499.75
$ws.Cells.Item(133, 10).Value = 68000
$ws.Cells.Item(133, 12).Value = 68000
$ws.Cells.Item(133, 14).Value = -73060

# CRP!row136
$ws.Cells.Item(136, 8).Value = 3691.4285
$ws.Cells.Item(136, 9).Value = 2726.9473
$ws.Cells.Item(136, 11).Value = 8180.841899999999
$ws.Cells.Item(136, 13).Value = -5630.841899999999

$ws = $wb.Worksheets.Item("CUL")
# CUL!row22
$ws.Cells.Item(22, 8).Value = 811.5
$ws.Cells.Item(22, 10).Value = 811.5
$ws.Cells.Item(22, 12).Value = 2434.5
$ws.Cells.Item(22, 14).Value = -2772.5

# CUL!row27
$ws.Cells.Item(27, 8).Value = 811.5
$ws.Cells.Item(27, 10).Value = 811.5
$ws.Cells.Item(27, 12).Value = 2434.5
$ws.Cells.Item(27, 14).Value = -2638.5

$ws = $wb.Worksheets.Item("GSM")
# GSM!row80
$ws.Cells.Item(80, 8).Value = 4719.75
$ws.Cells.Item(80, 10).Value = 4979.3335
$ws.Cells.Item(80, 12).Value = 4979.3335
$ws.Cells.Item(80, 14).Value = -6975.3335

# GSM!row83
$ws.Cells.Item(83, 8).Value = 4719.75
$ws.Cells.Item(83, 10).Value = 4979.3335
$ws.Cells.Item(83, 12).Value = 24896.6675
$ws.Cells.Item(83, 14).Value = -34880.6675

# GSM!row102
$ws.Cells.Item(102, 8).Value = 1335.6666
$ws.Cells.Item(102, 9).Value = 1231.0714
$ws.Cells.Item(102, 11).Value = 1231.0714
$ws.Cells.Item(102, 13).Value = 390.9286

# GSM!row107
$ws.Cells.Item(107, 8).Value = 796.2
$ws.Cells.Item(107, 9).Value = 882
$ws.Cells.Item(107, 10).Value = 667.5
$ws.Cells.Item(107, 11).Value = 882
$ws.Cells.Item(107, 12).Value = 667.5
$ws.Cells.Item(107, 13).Value = 1038
$ws.Cells.Item(107, 14).Value = -4507.5

# GSM!row126
$ws.Cells.Item(126, 8).Value = 2957.3125
$ws.Cells.Item(126, 9).Value = 2843.818
$ws.Cells.Item(126, 11).Value = 8531.454000000002
$ws.Cells.Item(126, 13).Value = -6061.454000000002

# GSM!row132
$ws.Cells.Item(132, 8).Value = 5189.5
$ws.Cells.Item(132, 9).Value = 2874.913
$ws.Cells.Item(132, 10).Value = 9284.538
$ws.Cells.Item(132, 11).Value = 8624.739
$ws.Cells.Item(132, 12).Value = 27853.614
$ws.Cells.Item(132, 13).Value = -6094.739
$ws.Cells.Item(132, 14).Value = -32913.614

$ws = $wb.Worksheets.Item("LTW")
# LTW!row35
$ws.Cells.Item(35, 8).Value = 1419.3334
$ws.Cells.Item(35, 9).Value = 1503.4
$ws.Cells.Item(35, 11).Value = 1503.4
$ws.Cells.Item(35, 13).Value = -1167.4

# LTW!row48
$ws.Cells.Item(48, 8).Value = 12249.5
$ws.Cells.Item(48, 9).Value = 12249.5
$ws.Cells.Item(48, 11).Value = 12249.5
$ws.Cells.Item(48, 13).Value = -11588.5

# LTW!row109
$ws.Cells.Item(109, 8).Value = 83399.5
$ws.Cells.Item(109, 10).Value = 83399.5
$ws.Cells.Item(109, 12).Value = 83399.5
$ws.Cells.Item(109, 14).Value = -86173.5

# LTW!row132
$ws.Cells.Item(132, 8).Value = 7112.5
$ws.Cells.Item(132, 9).Value = 2073.8462
$ws.Cells.Item(132, 10).Value = 20213
$ws.Cells.Item(132, 11).Value = 6221.5386
$ws.Cells.Item(132, 12).Value = 60639
$ws.Cells.Item(132, 13).Value = -3691.5386
$ws.Cells.Item(132, 14).Value = -65699

# LTW!row133
$ws.Cells.Item(133, 8).Value = 113055
$ws.Cells.Item(133, 10).Value = 113055
$ws.Cells.Item(133, 12).Value = 113055
$ws.Cells.Item(133, 14).Value = -118115

$ws = $wb.Worksheets.Item("WVR")
# WVR!row14
$ws.Cells.Item(14, 8).Value = 10504
$ws.Cells.Item(14, 9).Value = 10504
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 10504
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -10336
$ws.Cells.Item(14, 14).ClearContents()

# WVR!row107
$ws.Cells.Item(107, 8).Value = 1720.9166
$ws.Cells.Item(107, 9).Value = 909.25
$ws.Cells.Item(107, 11).Value = 2727.75
$ws.Cells.Item(107, 13).Value = -807.75

# WVR!row113
$ws.Cells.Item(113, 8).Value = 845.25
$ws.Cells.Item(113, 10).Value = 1407.375
$ws.Cells.Item(113, 12).Value = 4222.125
$ws.Cells.Item(113, 14).Value = -8562.125

# WVR!row126
$ws.Cells.Item(126, 8).Value = 3627.4443
$ws.Cells.Item(126, 9).Value = 3705.875
$ws.Cells.Item(126, 11).Value = 11117.625
$ws.Cells.Item(126, 13).Value = -8647.625

# WVR!row132
$ws.Cells.Item(132, 8).Value = 5545.5864
$ws.Cells.Item(132, 9).Value = 4682.625
$ws.Cells.Item(132, 11).Value = 14047.875
$ws.Cells.Item(132, 13).Value = -11517.875

# WVR!row136
$ws.Cells.Item(136, 8).Value = 2343.3428
$ws.Cells.Item(136, 9).Value = 722.0357
$ws.Cells.Item(136, 11).Value = 2166.1071
$ws.Cells.Item(136, 13).Value = 383.8928999999998
